# Slide 3 ("Wymagania funkcjonalne Aplikacji") content placeholder:
# every bullet line's run gets its Latin typeface pinned to "Aptos"
# (the new Office default font), matching what PowerPoint does when a
# user re-touches the text (select-all + reapply the theme body font).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$paraCount = $tr.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Length -gt 0) {
        $para.Font.Name = "Aptos"
    }
}
